$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending order of "Periodo Mora" values (previously descending 2003..1607,
# now reordered ascending 1607..2003) for rows 16 through 60 in column E.
$periodos = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$startRow = 16
for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
}
